$d = $word.ActiveDocument

# 1. Fix the "broken english": "...dashCommerce can has the following achievements..."
#    should read "...dashCommerce has the following achievements...".
$fix = $d.Content
$fix.Find.Execute("dashCommerce can has the following", $true, $false, $false, $false, $false, $true, 1, $false, "dashCommerce has the following", 2)

# 2. Re-split the run that now spans "...you might ask? Well, dashCommerce has the
#    following..." so it matches the author's updated run layout (break happens
#    mid-word, right after "you mi"). Nudging Font.Size away and back forces Word
#    to break the run at that point without altering the effective formatting.
$find = $d.Content
$find.Find.Execute("you might ask", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $find.Start + 6
$beforeSplit = $d.Range(0, $splitPos)
$beforeSplit.Font.Size = 99
$beforeSplit.Font.Size = 9

# 3. Move the "_GoBack" bookmark so it now sits right before "has the following
#    achievements to its credit:" (it is a singleton bookmark, so adding it here
#    removes it from its old spot near "...creating your online store!" for us).
$hasRange = $d.Content
$hasRange.Find.Execute("has the following achievements", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bookmarkSpot = $d.Range($hasRange.Start, $hasRange.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot)
